# "updated calib sheet and elec calibration"
#
# 1. The EDLE sheet (electricity calibration) value for "all electricity
#    sources" logit exponent is recalibrated from -3 to 1.
# 2. The EDLE sheet becomes the active/selected sheet (the calibration
#    sheet the author was last working in), moving the selection away
#    from the "About" sheet.

$wb = $excel.ActiveWorkbook

$wsEdle = $wb.Worksheets.Item("EDLE")

# Update the calibrated logit exponent value.
$wsEdle.Range("B2").Value = 1

# Make EDLE the active sheet/tab (was previously "About").
$wsEdle.Activate()
$wsEdle.Range("B3").Select()
